$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added at the top of this price-history
# block (row 183), pushing the existing rows 183-213 down by one
# (to 184-214). Excel's native row Insert keeps everything else
# (values, shared-string text, date styling, dimension) intact and
# shifted automatically, so we only need to populate the freshly
# inserted row with its own data.
$ws.Rows("183").Insert()

$row = 183
$ws.Cells.Item($row, 1).Value  = 5
$ws.Cells.Item($row, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item($row, 3).Value  = "Maule"
$ws.Cells.Item($row, 4).Value  = 44617
$ws.Cells.Item($row, 5).Value  = 7
$ws.Cells.Item($row, 6).Value  = "Fruta"
$ws.Cells.Item($row, 7).Value  = 100108
$ws.Cells.Item($row, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item($row, 9).Value  = 100108005
$ws.Cells.Item($row, 10).Value = "Piña"
$ws.Cells.Item($row, 11).Value = "Caramelo"
$ws.Cells.Item($row, 12).Value = "Segunda"
$ws.Cells.Item($row, 13).Value = 210
$ws.Cells.Item($row, 14).Value = 16000
$ws.Cells.Item($row, 15).Value = 16000
$ws.Cells.Item($row, 16).Value = 16000
$ws.Cells.Item($row, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item($row, 18).Value = "Ecuador"
$ws.Cells.Item($row, 19).Value = 1143
$ws.Cells.Item($row, 20).Value = 14
